# Rebuild the LR-pair data table (rows 2-21) with the refreshed TPM-based values.
# The sending/target cluster combinations are reordered and a new "MuSCs" target
# row is added for every sending cluster, extending the sheet from 16 to 21 rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> FAPs
$ws.Range("A2").Value = 'ECs'
$ws.Range("B2").Value = 'Icam1'
$ws.Range("C2").Value = 'Itgam'
$ws.Range("D2").Value = 'FAPs'
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 35.160799
$ws.Range("H2").Value = 105.482397
$ws.Range("I2").Value = 0.2238945559395223
$ws.Range("J2").Value = 0.2238945559395223
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.142723
$ws.Range("N2").Value = 0.428169
$ws.Range("O2").Value = 0.0009642800942465787
$ws.Range("P2").Value = 0.0009642800942465787
$ws.Range("Q2").Value = 5.018254715677
$ws.Range("R2").Value = 45.164292441093
$ws.Range("S2").Value = 0.0002158970635026584
$ws.Range("T2").Value = 0.0002158970635026584

# Row 3: ECs -> Inflammatory-Mac
$ws.Range("A3").Value = 'ECs'
$ws.Range("B3").Value = 'Icam1'
$ws.Range("C3").Value = 'Itgam'
$ws.Range("D3").Value = 'Inflammatory-Mac'
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 35.160799
$ws.Range("H3").Value = 105.482397
$ws.Range("I3").Value = 0.2238945559395223
$ws.Range("J3").Value = 0.2238945559395223
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 86.42780700000002
$ws.Range("N3").Value = 259.283421
$ws.Range("O3").Value = 0.5839326098770704
$ws.Range("P3").Value = 0.5839326098770704
$ws.Range("Q3").Value = 3038.870749937793
$ws.Range("R3").Value = 27349.83674944014
$ws.Range("S3").Value = 0.130739332387033
$ws.Range("T3").Value = 0.130739332387033

# Row 4: ECs -> MuSCs
$ws.Range("A4").Value = 'ECs'
$ws.Range("B4").Value = 'Icam1'
$ws.Range("C4").Value = 'Itgam'
$ws.Range("D4").Value = 'MuSCs'
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 35.160799
$ws.Range("H4").Value = 105.482397
$ws.Range("I4").Value = 0.2238945559395223
$ws.Range("J4").Value = 0.2238945559395223
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.006361333333333333
$ws.Range("N4").Value = 0.019084
$ws.Range("O4").Value = 0.00004297910712499435
$ws.Range("P4").Value = 0.00004297910712499435
$ws.Range("Q4").Value = 0.2236695627053333
$ws.Range("R4").Value = 2.013026064348
$ws.Range("S4").Value = 0.000009622788104427767
$ws.Range("T4").Value = 0.000009622788104427767

# Row 5: ECs -> Resolving-Mac
$ws.Range("A5").Value = 'ECs'
$ws.Range("B5").Value = 'Icam1'
$ws.Range("C5").Value = 'Itgam'
$ws.Range("D5").Value = 'Resolving-Mac'
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 35.160799
$ws.Range("H5").Value = 105.482397
$ws.Range("I5").Value = 0.2238945559395223
$ws.Range("J5").Value = 0.2238945559395223
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 61.43300833333333
$ws.Range("N5").Value = 184.299025
$ws.Range("O5").Value = 0.415060130921558
$ws.Range("P5").Value = 0.415060130921558
$ws.Range("Q5").Value = 2160.033657973658
$ws.Range("R5").Value = 19440.30292176292
$ws.Range("S5").Value = 0.09292970370088219
$ws.Range("T5").Value = 0.09292970370088219

# Row 6: FAPs -> FAPs
$ws.Range("A6").Value = 'FAPs'
$ws.Range("B6").Value = 'Icam1'
$ws.Range("C6").Value = 'Itgam'
$ws.Range("D6").Value = 'FAPs'
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 32.208719
$ws.Range("H6").Value = 96.626157
$ws.Range("I6").Value = 0.2050965007332699
$ws.Range("J6").Value = 0.2050965007332699
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.142723
$ws.Range("N6").Value = 0.428169
$ws.Range("O6").Value = 0.0009642800942465787
$ws.Range("P6").Value = 0.0009642800942465787
$ws.Range("Q6").Value = 4.596925001837
$ws.Range("R6").Value = 41.372325016533
$ws.Range("S6").Value = 0.000197770473056721
$ws.Range("T6").Value = 0.000197770473056721

# Row 7: FAPs -> Inflammatory-Mac
$ws.Range("A7").Value = 'FAPs'
$ws.Range("B7").Value = 'Icam1'
$ws.Range("C7").Value = 'Itgam'
$ws.Range("D7").Value = 'Inflammatory-Mac'
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 32.208719
$ws.Range("H7").Value = 96.626157
$ws.Range("I7").Value = 0.2050965007332699
$ws.Range("J7").Value = 0.2050965007332699
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 86.42780700000002
$ws.Range("N7").Value = 259.283421
$ws.Range("O7").Value = 0.5839326098770704
$ws.Range("P7").Value = 0.5839326098770704
$ws.Range("Q7").Value = 2783.728949449234
$ws.Range("R7").Value = 25053.5605450431
$ws.Range("S7").Value = 0.1197625349498328
$ws.Range("T7").Value = 0.1197625349498328

# Row 8: FAPs -> MuSCs
$ws.Range("A8").Value = 'FAPs'
$ws.Range("B8").Value = 'Icam1'
$ws.Range("C8").Value = 'Itgam'
$ws.Range("D8").Value = 'MuSCs'
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 32.208719
$ws.Range("H8").Value = 96.626157
$ws.Range("I8").Value = 0.2050965007332699
$ws.Range("J8").Value = 0.2050965007332699
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.006361333333333333
$ws.Range("N8").Value = 0.019084
$ws.Range("O8").Value = 0.00004297910712499435
$ws.Range("P8").Value = 0.00004297910712499435
$ws.Range("Q8").Value = 0.2048903977986667
$ws.Range("R8").Value = 1.844013580188
$ws.Range("S8").Value = 0.00000881486447597669
$ws.Range("T8").Value = 0.00000881486447597669

# Row 9: FAPs -> Resolving-Mac
$ws.Range("A9").Value = 'FAPs'
$ws.Range("B9").Value = 'Icam1'
$ws.Range("C9").Value = 'Itgam'
$ws.Range("D9").Value = 'Resolving-Mac'
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 32.208719
$ws.Range("H9").Value = 96.626157
$ws.Range("I9").Value = 0.2050965007332699
$ws.Range("J9").Value = 0.2050965007332699
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 61.43300833333333
$ws.Range("N9").Value = 184.299025
$ws.Range("O9").Value = 0.415060130921558
$ws.Range("P9").Value = 0.415060130921558
$ws.Range("Q9").Value = 1978.678502732992
$ws.Range("R9").Value = 17808.10652459693
$ws.Range("S9").Value = 0.0851273804459044
$ws.Range("T9").Value = 0.0851273804459044

# Row 10: Inflammatory-Mac -> FAPs
$ws.Range("A10").Value = 'Inflammatory-Mac'
$ws.Range("B10").Value = 'Icam1'
$ws.Range("C10").Value = 'Itgam'
$ws.Range("D10").Value = 'FAPs'
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 51.53356533333334
$ws.Range("H10").Value = 154.600696
$ws.Range("I10").Value = 0.3281519491717758
$ws.Range("J10").Value = 0.3281519491717758
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.142723
$ws.Range("N10").Value = 0.428169
$ws.Range("O10").Value = 0.0009642800942465787
$ws.Range("P10").Value = 0.0009642800942465787
$ws.Range("Q10").Value = 7.355025045069334
$ws.Range("R10").Value = 66.195225405624
$ws.Range("S10").Value = 0.0003164303924745585
$ws.Range("T10").Value = 0.0003164303924745585

# Row 11: Inflammatory-Mac -> Inflammatory-Mac
$ws.Range("A11").Value = 'Inflammatory-Mac'
$ws.Range("B11").Value = 'Icam1'
$ws.Range("C11").Value = 'Itgam'
$ws.Range("D11").Value = 'Inflammatory-Mac'
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 51.53356533333334
$ws.Range("H11").Value = 154.600696
$ws.Range("I11").Value = 0.3281519491717758
$ws.Range("J11").Value = 0.3281519491717758
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 86.42780700000002
$ws.Range("N11").Value = 259.283421
$ws.Range("O11").Value = 0.5839326098770704
$ws.Range("P11").Value = 0.5839326098770704
$ws.Range("Q11").Value = 4453.933038651225
$ws.Range("R11").Value = 40085.39734786102
$ws.Range("S11").Value = 0.1916186241161228
$ws.Range("T11").Value = 0.1916186241161228

# Row 12: Inflammatory-Mac -> MuSCs
$ws.Range("A12").Value = 'Inflammatory-Mac'
$ws.Range("B12").Value = 'Icam1'
$ws.Range("C12").Value = 'Itgam'
$ws.Range("D12").Value = 'MuSCs'
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 51.53356533333334
$ws.Range("H12").Value = 154.600696
$ws.Range("I12").Value = 0.3281519491717758
$ws.Range("J12").Value = 0.3281519491717758
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.006361333333333333
$ws.Range("N12").Value = 0.019084
$ws.Range("O12").Value = 0.00004297910712499435
$ws.Range("P12").Value = 0.00004297910712499435
$ws.Range("Q12").Value = 0.3278221869404445
$ws.Range("R12").Value = 2.950399682464
$ws.Range("S12").Value = 0.00001410367777672945
$ws.Range("T12").Value = 0.00001410367777672945

# Row 13: Inflammatory-Mac -> Resolving-Mac
$ws.Range("A13").Value = 'Inflammatory-Mac'
$ws.Range("B13").Value = 'Icam1'
$ws.Range("C13").Value = 'Itgam'
$ws.Range("D13").Value = 'Resolving-Mac'
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 51.53356533333334
$ws.Range("H13").Value = 154.600696
$ws.Range("I13").Value = 0.3281519491717758
$ws.Range("J13").Value = 0.3281519491717758
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 61.43300833333333
$ws.Range("N13").Value = 184.299025
$ws.Range("O13").Value = 0.415060130921558
$ws.Range("P13").Value = 0.415060130921558
$ws.Range("Q13").Value = 3165.861948569045
$ws.Range("R13").Value = 28492.7575371214
$ws.Range("S13").Value = 0.1362027909854017
$ws.Range("T13").Value = 0.1362027909854017

# Row 14: MuSCs -> FAPs
$ws.Range("A14").Value = 'MuSCs'
$ws.Range("B14").Value = 'Icam1'
$ws.Range("C14").Value = 'Itgam'
$ws.Range("D14").Value = 'FAPs'
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.5955593333333333
$ws.Range("H14").Value = 1.786678
$ws.Range("I14").Value = 0.003792362411113143
$ws.Range("J14").Value = 0.003792362411113143
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.142723
$ws.Range("N14").Value = 0.428169
$ws.Range("O14").Value = 0.0009642800942465787
$ws.Range("P14").Value = 0.0009642800942465787
$ws.Range("Q14").Value = 0.08500001473133334
$ws.Range("R14").Value = 0.7650001325820001
$ws.Range("S14").Value = 0.000003656899583205364
$ws.Range("T14").Value = 0.000003656899583205364

# Row 15: MuSCs -> Inflammatory-Mac
$ws.Range("A15").Value = 'MuSCs'
$ws.Range("B15").Value = 'Icam1'
$ws.Range("C15").Value = 'Itgam'
$ws.Range("D15").Value = 'Inflammatory-Mac'
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.5955593333333333
$ws.Range("H15").Value = 1.786678
$ws.Range("I15").Value = 0.003792362411113143
$ws.Range("J15").Value = 0.003792362411113143
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 86.42780700000002
$ws.Range("N15").Value = 259.283421
$ws.Range("O15").Value = 0.5839326098770704
$ws.Range("P15").Value = 0.5839326098770704
$ws.Range("Q15").Value = 51.47288711838201
$ws.Range("R15").Value = 463.255984065438
$ws.Range("S15").Value = 0.002214484080320997
$ws.Range("T15").Value = 0.002214484080320997

# Row 16: MuSCs -> MuSCs
$ws.Range("A16").Value = 'MuSCs'
$ws.Range("B16").Value = 'Icam1'
$ws.Range("C16").Value = 'Itgam'
$ws.Range("D16").Value = 'MuSCs'
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.5955593333333333
$ws.Range("H16").Value = 1.786678
$ws.Range("I16").Value = 0.003792362411113143
$ws.Range("J16").Value = 0.003792362411113143
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.006361333333333333
$ws.Range("N16").Value = 0.019084
$ws.Range("O16").Value = 0.00004297910712499435
$ws.Range("P16").Value = 0.00004297910712499435
$ws.Range("Q16").Value = 0.003788551439111111
$ws.Range("R16").Value = 0.034096962952
$ws.Range("S16").Value = 0.0000001629923503240337
$ws.Range("T16").Value = 0.0000001629923503240337

# Row 17: MuSCs -> Resolving-Mac
$ws.Range("A17").Value = 'MuSCs'
$ws.Range("B17").Value = 'Icam1'
$ws.Range("C17").Value = 'Itgam'
$ws.Range("D17").Value = 'Resolving-Mac'
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.5955593333333333
$ws.Range("H17").Value = 1.786678
$ws.Range("I17").Value = 0.003792362411113143
$ws.Range("J17").Value = 0.003792362411113143
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 61.43300833333333
$ws.Range("N17").Value = 184.299025
$ws.Range("O17").Value = 0.415060130921558
$ws.Range("P17").Value = 0.415060130921558
$ws.Range("Q17").Value = 36.58700148766111
$ws.Range("R17").Value = 329.28301338895
$ws.Range("S17").Value = 0.001574058438858616
$ws.Range("T17").Value = 0.001574058438858616

# Row 18: Resolving-Mac -> FAPs
$ws.Range("A18").Value = 'Resolving-Mac'
$ws.Range("B18").Value = 'Icam1'
$ws.Range("C18").Value = 'Itgam'
$ws.Range("D18").Value = 'FAPs'
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 37.54313466666667
$ws.Range("H18").Value = 112.629404
$ws.Range("I18").Value = 0.2390646317443189
$ws.Range("J18").Value = 0.2390646317443189
$ws.Range("K18").Value = 1
$ws.Range("L18").Value = 0.3333333333333333
$ws.Range("M18").Value = 0.142723
$ws.Range("N18").Value = 0.428169
$ws.Range("O18").Value = 0.0009642800942465787
$ws.Range("P18").Value = 0.0009642800942465787
$ws.Range("Q18").Value = 5.358268809030667
$ws.Range("R18").Value = 48.224419281276
$ws.Range("S18").Value = 0.0002305252656294355
$ws.Range("T18").Value = 0.0002305252656294355

# Row 19: Resolving-Mac -> Inflammatory-Mac
$ws.Range("A19").Value = 'Resolving-Mac'
$ws.Range("B19").Value = 'Icam1'
$ws.Range("C19").Value = 'Itgam'
$ws.Range("D19").Value = 'Inflammatory-Mac'
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 37.54313466666667
$ws.Range("H19").Value = 112.629404
$ws.Range("I19").Value = 0.2390646317443189
$ws.Range("J19").Value = 0.2390646317443189
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 86.42780700000002
$ws.Range("N19").Value = 259.283421
$ws.Range("O19").Value = 0.5839326098770704
$ws.Range("P19").Value = 0.5839326098770704
$ws.Range("Q19").Value = 3244.770797145677
$ws.Range("R19").Value = 29202.93717431109
$ws.Range("S19").Value = 0.1395976343437609
$ws.Range("T19").Value = 0.1395976343437609

# Row 20: Resolving-Mac -> MuSCs
$ws.Range("A20").Value = 'Resolving-Mac'
$ws.Range("B20").Value = 'Icam1'
$ws.Range("C20").Value = 'Itgam'
$ws.Range("D20").Value = 'MuSCs'
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 37.54313466666667
$ws.Range("H20").Value = 112.629404
$ws.Range("I20").Value = 0.2390646317443189
$ws.Range("J20").Value = 0.2390646317443189
$ws.Range("K20").Value = 1
$ws.Range("L20").Value = 0.3333333333333333
$ws.Range("M20").Value = 0.006361333333333333
$ws.Range("N20").Value = 0.019084
$ws.Range("O20").Value = 0.00004297910712499435
$ws.Range("P20").Value = 0.00004297910712499435
$ws.Range("Q20").Value = 0.2388243939928889
$ws.Range("R20").Value = 2.149419545936
$ws.Range("S20").Value = 0.00001027478441753641
$ws.Range("T20").Value = 0.00001027478441753641

# Row 21: Resolving-Mac -> Resolving-Mac
$ws.Range("A21").Value = 'Resolving-Mac'
$ws.Range("B21").Value = 'Icam1'
$ws.Range("C21").Value = 'Itgam'
$ws.Range("D21").Value = 'Resolving-Mac'
$ws.Range("E21").Value = 3
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 37.54313466666667
$ws.Range("H21").Value = 112.629404
$ws.Range("I21").Value = 0.2390646317443189
$ws.Range("J21").Value = 0.2390646317443189
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 61.43300833333333
$ws.Range("N21").Value = 184.299025
$ws.Range("O21").Value = 0.415060130921558
$ws.Range("P21").Value = 0.415060130921558
$ws.Range("Q21").Value = 2306.387704836789
$ws.Range("R21").Value = 20757.4893435311
$ws.Range("S21").Value = 0.09922619735051107
$ws.Range("T21").Value = 0.09922619735051107
